# Generate Report for Handback
# The "9afdfc77-1c64-4270-878e-ca47b051a6ea.md" file has moved from
# "Ready for handoff" to "Handed back: in sync with en-US". Update its
# Status and Latest Handback DateTime columns on every worksheet.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Range("H3").Value = "2016-03-17 12:35:50"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Range("H3").Value = "2016-03-17 12:35:56"
